$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 411, shifting existing rows 411:514 down to 412:515
$ws.Rows.Item(411).Insert()

# Populate the newly inserted row 411 with its data
$ws.Cells.Item(411, 1).Value = 5
$ws.Cells.Item(411, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(411, 3).Value = "Maule"
$ws.Cells.Item(411, 4).Value = 45204
$ws.Cells.Item(411, 5).Value = 7
$ws.Cells.Item(411, 6).Value = 100112045
$ws.Cells.Item(411, 7).Value = "Zapallo"
$ws.Cells.Item(411, 8).Value = "Paine"
$ws.Cells.Item(411, 9).Value = "1a (guarda)"
$ws.Cells.Item(411, 10).Value = 1500
$ws.Cells.Item(411, 11).Value = 500
$ws.Cells.Item(411, 12).Value = 500
$ws.Cells.Item(411, 13).Value = 500
$ws.Cells.Item(411, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(411, 15).Value = "Región del Maule"
$ws.Cells.Item(411, 16).Value = 500
$ws.Cells.Item(411, 17).Value = 1
$ws.Cells.Item(411, 18).Value = "Hortaliza"
